$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cells in column D whose refreshed value looks like a plain number must be forced to
# Text format before the assignment, otherwise Excel auto-converts the literal string
# into a floating point number (losing the original fixed 2-decimal text formatting).

$ws.Range("D2").Value = "58.253.33"
$ws.Range("D3").Value = "2.610.95"
$ws.Range("E3").Value = "  -3.38%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.84"
$ws.Range("E5").Value = "  -1.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.91"
$ws.Range("E6").Value = "  -1.77%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("E8").Value = "  -1.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.68"
$ws.Range("E9").Value = "  -1.62%  "
$ws.Range("E10").Value = "  -2.39%  "
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("D13").Value = "3.067.95"
$ws.Range("E13").Value = "  -3.62%  "
$ws.Range("D14").Value = "58.224.12"
$ws.Range("E14").Value = "  -3.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.99"
$ws.Range("E15").Value = "  -1.25%  "
$ws.Range("E16").Value = "  -1.25%  "
$ws.Range("D17").Value = "2.607.40"
$ws.Range("E17").Value = "  -4.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.41"
$ws.Range("E18").Value = "  -2.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "335.82"
$ws.Range("E19").Value = "  -2.47%  "
$ws.Range("E20").Value = "  -2.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.27"
$ws.Range("E21").Value = "  -3.17%  "
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.41"
$ws.Range("E23").Value = "  +1.70%  "
$ws.Range("E24").Value = "  -1.26%  "
$ws.Range("E25").Value = "  -1.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.996"
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.12"
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("D28").Value = "0.0₃0792"
$ws.Range("E28").Value = "  -3.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.65"
$ws.Range("E29").Value = "  -2.38%  "
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("E31").Value = "  -1.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.77"
$ws.Range("E32").Value = "  -1.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "150.67"
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.10"
$ws.Range("E34").Value = "  -3.72%  "
$ws.Range("E35").Value = "  -3.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.889"
$ws.Range("E36").Value = "  -5.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.851"
$ws.Range("E37").Value = "  -2.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.28"
$ws.Range("E38").Value = "  -2.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.43"
$ws.Range("E39").Value = "  -5.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.63"
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("E42").Value = "  -1.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0966"
$ws.Range("E43").Value = "  -2.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "269.31"
$ws.Range("E44").Value = "  -3.95%  "
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.15"
$ws.Range("E46").Value = "  -4.96%  "
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("D48").Value = "2.037.17"
$ws.Range("E48").Value = "  -4.98%  "
$ws.Range("E49").Value = "  -1.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.68"
$ws.Range("E50").Value = "  -3.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.24"
$ws.Range("E51").Value = "  -4.31%  "
